$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Cells.Item(52, 1).Value = "Linking_AutoUser"
$ws.Cells.Item(52, 2).Value = "Password1"
$ws.Cells.Item(52, 5).Value = "Default user for Linking tests"
$ws.Cells.Item(52, 6).Value = "N"
$ws.Cells.Item(52, 7).Value = "linking.autouser@mailinator.com"
